$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267 (shifts existing rows 267..278 down to 268..279)
$ws.Rows.Item(267).Insert()

# Populate the new row 267 with the weekly price-report record
$ws.Range("A267").Value = 9
$ws.Range("B267").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C267").Value = "Metropolitana"
$ws.Range("D267").Value = 44714
$ws.Range("E267").Value = 13
$ws.Range("F267").Value = 100112001
$ws.Range("G267").Value = "Berenjena"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 70
$ws.Range("K267").Value = 9000
$ws.Range("L267").Value = 10000
$ws.Range("M267").Value = 9500
$ws.Range("N267").Value = "`$/caja 50 unidades"
$ws.Range("O267").Value = "Región de Arica y Parinacota"
$ws.Range("P267").Value = 190
$ws.Range("Q267").Value = 50
$ws.Range("R267").Value = "Hortaliza"
